$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3 / Row 4: Procedure.performedPeriod -> split into .start / .end ---
$ws.Rows(3).RowHeight = 25.5
$ws.Range("J3").Value = "Procedure.performedPeriod.start"
$ws.Range("K3").Value = ""

$ws.Rows(4).RowHeight = 25.5
$ws.Range("J4").Value = "Procedure.performedPeriod.end"
$ws.Range("K4").Value = ""

# --- Rows 5-10: clear stale "Fits" / ZIB reference notes (text unchanged in J col) ---
$ws.Range("K5").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("K7").Value = "Current valueset misses valid OID for system. ZIBs 3.0 has fixed the issue. "
$ws.Range("K8").Value = ""
$ws.Range("K9").Value = ""
$ws.Range("K10").Value = ""

# --- Row 11: update note text ---
$ws.Range("K11").Value = "Still need an extension  / Or will procedureRequest be sufficient?"

# --- Old leftover notes rows 13-16: wipe content+format, then reuse formatting ---
# that already exists further down the sheet (avoids creating brand-new style
# records, mirroring what Excel does when you paste formats from a nearby cell).
$ws.Range("K13:K16").Clear()

$ws.Range("A33").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("B32").Copy()
$ws.Range("B13").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Delete the duplicated second copy of the table (rows 19-35) ---
$ws.Rows("19:35").Delete()

# --- Column J width adjustment ---
$ws.Columns(10).ColumnWidth = 30.28515625
